$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UserLog")

$data = @(
    @('E20241119701802024002172', 'Paul Atreides', 'E20047030021642A023C3D52', 'Vernier Caliper', 'BORROW', '-', '2025-11-07 08:44:57'),
    @('E20241119701802024002172', 'Paul Atreides', 'E2004702FFA1642A023C3D4A', 'Micrometer', 'BORROW', '-', '2025-11-07 08:45:05'),
    @('E20241119701802024002172', 'Paul Atreides', 'E20047030021642A023C3D52', 'Vernier Caliper', 'RETURN', 'Good', '2025-11-07 08:45:26'),
    @('E20241119701802024002172', 'Paul Atreides', 'E2004702FFA1642A023C3D4A', 'Micrometer', 'RETURN', 'Damaged', '2025-11-07 08:45:33'),
    @('E20047124741662E032E737C', 'Kevin Flynn', 'E2004702FFA1642A023C3D4A', 'Micrometer', 'BORROW', '-', '2025-11-07 08:47:33'),
    @('E20047124741662E032E737C', 'Kevin Flynn', 'E20047030021642A023C3D52', 'Vernier Caliper', 'BORROW', '-', '2025-11-07 08:47:34'),
    @('E20047124741662E032E737C', 'Kevin Flynn', 'E2004702FFA1642A023C3D4A', 'Micrometer', 'RETURN', 'Good', '2025-11-07 08:48:05'),
    @('E20047124741662E032E737C', 'Kevin Flynn', 'E20047030021642A023C3D52', 'Vernier Caliper', 'RETURN', 'Good', '2025-11-07 08:48:09'),
    @('E20241119701802024000520', 'Anakin Skywalker', 'E2004702FFA1642A023C3D4A', 'Micrometer', 'BORROW', '-', '2025-11-07 08:54:32'),
    @('E20241119701802024000520', 'Anakin Skywalker', 'E20047030021642A023C3D52', 'Vernier Caliper', 'BORROW', '-', '2025-11-07 08:54:45'),
    @('E20241119701802024000520', 'Anakin Skywalker', 'E20047030021642A023C3D52', 'Vernier Caliper', 'RETURN', 'Good', '2025-11-07 08:55:21'),
    @('E20241119701802024000520', 'Anakin Skywalker', 'E2004702FFA1642A023C3D4A', 'Micrometer', 'RETURN', 'Good', '2025-11-07 08:55:30'),
    @('E20241119701802024000520', 'Anakin Skywalker', 'E20047030021642A023C3D52', 'Vernier Caliper', 'BORROW', '-', '2025-11-07 08:56:34'),
    @('E20241119701802024000520', 'Anakin Skywalker', 'E2004702FFA1642A023C3D4A', 'Micrometer', 'BORROW', '-', '2025-11-07 08:56:36'),
    @('E20241119701802024000520', 'Anakin Skywalker', 'E2004702FFA1642A023C3D4A', 'Micrometer', 'RETURN', 'Good', '2025-11-07 08:56:46'),
    @('E20241119701802024000520', 'Anakin Skywalker', 'E20047030021642A023C3D52', 'Vernier Caliper', 'RETURN', 'Good', '2025-11-07 08:56:55'),
    @('E20241119701802024000520', 'Anakin Skywalker', 'E20047030021642A023C3D52', 'Vernier Caliper', 'BORROW', '-', '2025-11-07 08:58:22'),
    @('E20241119701802024000520', 'Anakin Skywalker', 'E2004702FFA1642A023C3D4A', 'Micrometer', 'BORROW', '-', '2025-11-07 08:58:23'),
    @('E20241119701802024000520', 'Anakin Skywalker', 'E20047030021642A023C3D52', 'Vernier Caliper', 'RETURN', 'Good', '2025-11-07 08:58:32'),
    @('E20241119701802024000520', 'Anakin Skywalker', 'E2004702FFA1642A023C3D4A', 'Micrometer', 'RETURN', 'Damaged', '2025-11-07 08:58:40'),
    @('E20241119701802024000520', 'Anakin Skywalker', 'E2004702FFA1642A023C3D4A', 'Micrometer', 'BORROW', '-', '2025-11-07 09:00:45'),
    @('E20241119701802024000520', 'Anakin Skywalker', 'E20047030021642A023C3D52', 'Vernier Caliper', 'BORROW', '-', '2025-11-07 09:00:48'),
    @('E20241119701802024000520', 'Anakin Skywalker', 'E20047030021642A023C3D52', 'Vernier Caliper', 'RETURN', 'Good', '2025-11-07 09:00:54'),
    @('E20241119701802024000520', 'Anakin Skywalker', 'E2004702FFA1642A023C3D4A', 'Micrometer', 'RETURN', 'Good', '2025-11-07 09:01:19'),
    @('E20047124741662E032E737C', 'Kevin Flynn', 'E20047030021642A023C3D52', 'Vernier Caliper', 'BORROW', '-', '2025-11-07 09:04:04'),
    @('E20047124741662E032E737C', 'Kevin Flynn', 'E2004702FFA1642A023C3D4A', 'Micrometer', 'BORROW', '-', '2025-11-07 09:04:04'),
    @('E20047124741662E032E737C', 'Kevin Flynn', 'E2004702FFA1642A023C3D4A', 'Micrometer', 'RETURN', 'Good', '2025-11-07 09:04:14'),
    @('E20047124741662E032E737C', 'Kevin Flynn', 'E20047030021642A023C3D52', 'Vernier Caliper', 'RETURN', 'Good', '2025-11-07 09:04:20'),
    @('E20241119701802024000520', 'Anakin Skywalker', 'E2004702FFA1642A023C3D4A', 'Micrometer', 'BORROW', '-', '2025-11-07 09:08:43'),
    @('E20241119701802024000520', 'Anakin Skywalker', 'E20047030021642A023C3D52', 'Vernier Caliper', 'BORROW', '-', '2025-11-07 09:08:45'),
    @('E20241119701802024000520', 'Anakin Skywalker', 'E2004702FFA1642A023C3D4A', 'Micrometer', 'RETURN', 'Good', '2025-11-07 09:09:38'),
    @('E20241119701802024000520', 'Anakin Skywalker', 'E20047030021642A023C3D52', 'Vernier Caliper', 'RETURN', 'Damaged', '2025-11-07 09:09:43'),
    @('E20241119701802024002172', 'Paul Atreides', 'E2004702FFA1642A023C3D4A', 'Micrometer', 'BORROW', '-', '2025-11-07 09:09:51'),
    @('E20241119701802024002172', 'Paul Atreides', 'E20047030021642A023C3D52', 'Vernier Caliper', 'BORROW', '-', '2025-11-07 09:09:53'),
    @('E20241119701802024002172', 'Paul Atreides', 'E20047030021642A023C3D52', 'Vernier Caliper', 'RETURN', 'Good', '2025-11-07 09:10:28'),
    @('E20241119701802024002172', 'Paul Atreides', 'E2004702FFA1642A023C3D4A', 'Micrometer', 'RETURN', 'Good', '2025-11-07 09:12:58'),
    @('E20241119701802024000520', 'Anakin Skywalker', 'E2004702FFA1642A023C3D4A', 'Micrometer', 'BORROW', '-', '2025-11-07 09:13:11'),
    @('E20241119701802024000520', 'Anakin Skywalker', 'E20047030021642A023C3D52', 'Vernier Caliper', 'BORROW', '-', '2025-11-07 09:13:13'),
    @('E20241119701802024000520', 'Anakin Skywalker', 'E2004702FFA1642A023C3D4A', 'Micrometer', 'RETURN', 'Good', '2025-11-07 09:13:28'),
    @('E20241119701802024000520', 'Anakin Skywalker', 'E20047030021642A023C3D52', 'Vernier Caliper', 'RETURN', 'Good', '2025-11-07 09:13:42'),
    @('E20241119701802024000520', 'Anakin Skywalker', 'E20047030021642A023C3D52', 'Vernier Caliper', 'BORROW', '-', '2025-11-07 09:19:19'),
    @('E20241119701802024000520', 'Anakin Skywalker', 'E2004702FFA1642A023C3D4A', 'Micrometer', 'BORROW', '-', '2025-11-07 09:19:20'),
    @('E20241119701802024000520', 'Anakin Skywalker', 'E2004702FFA1642A023C3D4A', 'Micrometer', 'RETURN', 'Good', '2025-11-07 09:19:31'),
    @('E20241119701802024000520', 'Anakin Skywalker', 'E20047030021642A023C3D52', 'Vernier Caliper', 'RETURN', 'Good', '2025-11-07 09:19:34'),
)

$startRow = 28
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($col = 1; $col -le 7; $col++) {
        $ws.Cells.Item($row, $col).Value = $rowData[$col - 1]
    }
}
